# "implementação de excluir linhas zeradas e possibilidade de excluir clientes"
# - Remove the zeroed-out / removed debtor row (old row 5: "1", R$, 456).
# - Shift the remaining debtor rows up, with refreshed names/amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the now-unwanted last row outright (row count goes from 5 to 4).
$ws.Rows.Item(5).Delete()

function Set-DebtorRow([int]$row, [string]$name, [string]$currency, [string]$amount) {
    $ws.Range("A$row").Value = $name
    $ws.Range("B$row").Value = $currency
    # The amount column stores numeric-looking values as TEXT (matches the
    # original sheet's convention). Flip the cell to text first so Excel
    # doesn't auto-coerce the string into a Number, then drop the format
    # back to the default "Normal" style so no stray number-format/style
    # is left behind on the cell.
    $ws.Range("C$row").NumberFormat = "@"
    $ws.Range("C$row").Value = $amount
    $ws.Range("C$row").Style = "Normal"
}

Set-DebtorRow 1 "Joe"     "R$" "866.55"
Set-DebtorRow 2 "Jose"    "R$" "40"
Set-DebtorRow 3 "Arnaldo" "R$" "209"
Set-DebtorRow 4 "Maria"   "R$" "450"
